# Remove the "ACID Compliant" bullet point from the
# "Why MySQL over PostgreSQL or MongoDB?" list.
$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*ACID Compliant*MySQL ensures data integrity*") {
        $p.Range.Delete()
        break
    }
}
